# Update the cryptos list (Price / Volume(1h) columns) with refreshed data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  D = "27.214.20";  E = "  +0.89%  " },
    @{ Row = 3;  D = "1.688.39";   E = "  +0.43%  " },
    @{ Row = 4;  D = $null;        E = "  +0.01%  " },
    @{ Row = 5;  D = "215.65";     E = "  +0.21%  " },
    @{ Row = 6;  D = "0.520";      E = "  +0.40%  " },
    @{ Row = 7;  D = $null;        E = "  +0.03%  " },
    @{ Row = 8;  D = "23.16";      E = "  +10.61%  " },
    @{ Row = 9;  D = $null;        E = "  +4.64%  " },
    @{ Row = 10; D = $null;        E = "  +1.33%  " },
    @{ Row = 11; D = "0.0890";     E = "  +0.37%  " },
    @{ Row = 12; D = "1.926.56";   E = "  +0.45%  " },
    @{ Row = 13; D = "1.685.42";   E = "  +0.51%  " },
    @{ Row = 14; D = "4.21";       E = "  +2.45%  " },
    @{ Row = 15; D = "0.556";      E = "  +4.49%  " },
    @{ Row = 16; D = "67.12";      E = "  +1.90%  " },
    @{ Row = 17; D = "27.208.72";  E = "  +0.79%  " },
    @{ Row = 18; D = "237.06";     E = "  +0.36%  " },
    @{ Row = 19; D = $null;        E = "  -1.43%  " },
    @{ Row = 20; D = "0.0₃0747";   E = "  +1.68%  " },
    @{ Row = 21; D = $null;        E = "  -0.01%  " },
    @{ Row = 22; D = $null;        E = "  +3.56%  " },
    @{ Row = 23; D = "9.64";       E = "  +4.64%  " },
    @{ Row = 24; D = $null;        E = "  -1.66%  " },
    @{ Row = 25; D = "147.33";     E = "  +0.64%  " },
    @{ Row = 26; D = "7.32";       E = "  +1.36%  " },
    @{ Row = 27; D = "16.46";      E = "  +2.39%  " },
    @{ Row = 28; D = "0.113";      E = "  +0.58%  " },
    @{ Row = 29; D = "0.999";      E = "  -0.13%  " },
    @{ Row = 30; D = "0.0507";     E = "  +1.31%  " },
    @{ Row = 31; D = $null;        E = "  +0.39%  " },
    @{ Row = 32; D = $null;        E = "  +1.87%  " },
    @{ Row = 33; D = "1.551.73";   E = "  +3.85%  " },
    @{ Row = 34; D = $null;        E = "  +2.19%  " },
    @{ Row = 35; D = $null;        E = "  -0.89%  " },
    @{ Row = 36; D = "0.607";      E = "  +3.53%  " },
    @{ Row = 37; D = "0.950";      E = "  +3.52%  " },
    @{ Row = 38; D = $null;        E = "  -0.35%  " },
    @{ Row = 39; D = $null;        E = "  -0.66%  " },
    @{ Row = 40; D = $null;        E = "  +2.63%  " },
    @{ Row = 41; D = "69.44";      E = "  +2.87%  " },
    @{ Row = 42; D = "5.77";       E = "  +0.38%  " },
    @{ Row = 43; D = $null;        E = "  +0.03%  " },
    @{ Row = 44; D = $null;        E = "  -0.33%  " },
    @{ Row = 45; D = "1.835.38";   E = "  +0.78%  " },
    @{ Row = 46; D = $null;        E = "  +1.41%  " },
    @{ Row = 47; D = "90.83";      E = "  +0.23%  " },
    @{ Row = 48; D = "1.62";       E = "  +6.19%  " },
    @{ Row = 49; D = $null;        E = "  +5.28%  " },
    @{ Row = 50; D = "8.31";       E = "  +7.40%  " },
    @{ Row = 51; D = $null;        E = "  +1.54%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $dCell = $ws.Cells.Item($u.Row, 4)
        $dCell.NumberFormat = "@"
        $dCell.Value = $u.D
    }
    $eCell = $ws.Cells.Item($u.Row, 5)
    $eCell.NumberFormat = "@"
    $eCell.Value = $u.E
}
